$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.892999999999999
$ws.Range("C7").Value = -13.163
$ws.Range("E7").Value = 16.32
$ws.Range("A9").Value = -21.831
$ws.Range("E10").Value = 16.195
$ws.Range("C12").Value = -11.901
$ws.Range("E13").Value = 16.849
$ws.Range("C14").Value = -12.89
$ws.Range("D15").Value = -8.348000000000003
$ws.Range("E16").Value = 16.869
$ws.Range("A18").Value = -22.072
$ws.Range("A20").Value = -20.558
$ws.Range("E20").Value = 16.218
$ws.Range("E24").Value = 16.639
$ws.Range("C26").Value = -12.467
$ws.Range("A27").Value = -21.761
$ws.Range("C27").Value = -13.363
$ws.Range("C29").Value = -11.857
$ws.Range("D33").Value = -7.515000000000001
$ws.Range("A35").Value = -19.955
$ws.Range("D35").Value = -7.737
$ws.Range("C37").Value = -13.151
$ws.Range("C38").Value = -13.818
$ws.Range("D38").Value = -8.001999999999999
$ws.Range("E39").Value = 16.532
$ws.Range("D43").Value = -7.580999999999999
$ws.Range("D44").Value = -7.628
$ws.Range("D47").Value = -7.662999999999999
$ws.Range("E47").Value = 16.813
$ws.Range("E48").Value = 17.012
$ws.Range("C51").Value = -12.751
$ws.Range("D51").Value = -7.590000000000001
$ws.Range("C52").Value = -11.441
$ws.Range("E52").Value = 17.262
$ws.Range("C55").Value = -13.65
$ws.Range("E56").Value = 16.974
$ws.Range("D57").Value = -7.944999999999999
$ws.Range("D63").Value = -7.665000000000001
$ws.Range("A69").Value = -21.75
$ws.Range("C69").Value = -11.78
$ws.Range("C70").Value = -13.392
$ws.Range("D70").Value = -7.982000000000001
$ws.Range("A76").Value = -20.306
$ws.Range("A78").Value = -19.854
$ws.Range("C81").Value = -13.561
$ws.Range("A82").Value = -22.154
$ws.Range("A83").Value = -21.719
$ws.Range("C83").Value = -12.903
$ws.Range("E84").Value = 16.856
$ws.Range("D88").Value = -7.895
$ws.Range("A93").Value = -21.582
$ws.Range("D99").Value = -8.103999999999999
$ws.Range("E100").Value = 16.767
$ws.Range("E101").Value = 16.889
$ws.Range("C102").Value = -13.663
